{"js": "// Replace the huge Java stack-trace text that was pasted into the\n// \"aqlFeatureAccess(...) ... failed:\" run with the short message\n// \"Feature name not found in EClass EObject\", leaving the rest of the\n// document (and the run's bold/orange formatting, plus the trailing\n// <w:br/>) untouched.\n\nconst body = context.document.body;\n\n// Locate the start of the long text (the beginning of the offending run)\n// and the end of the long text (the last line of the stack trace, plus\n// the trailing newline character that precedes the run's closing tag).\nconst startResults = body.search(\n  \"aqlFeatureAccess(org.eclipse.emf.ecore.EObject,java.lang.String) with arguments\",\n  { matchCase: true }\n);\nstartResults.load(\"items\");\nawait context.sync();\n\nconst endResults = body.search(\"Main.main(Main.java:1472)\\n\", { matchCase: true });\nendResults.load(\"items\");\nawait context.sync();\n\nif (startResults.items.length === 0 || endResults.items.length === 0) {\n  throw new Error(\"Could not locate the stack-trace text to replace.\");\n}\n\nconst startRange = startResults.items[0];\nconst endRange = endResults.items[0];\n\n// Expand from the start of the stack trace all the way through to (and\n// including) the trailing newline right before the run ends, then swap\n// that whole span for the short replacement text (formatting of the\n// host run, i.e. bold + orange color, is preserved automatically).\nconst fullRange = startRange.expandTo(endRange);\nfullRange.insertText(\"Feature name not found in EClass EObject\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Replace the huge Java stack-trace text that was pasted into the\n# \"aqlFeatureAccess(...) ... failed:\" run with the short message\n# \"Feature name not found in EClass EObject\", leaving the rest of the\n# document (and the run's bold/orange formatting, plus the trailing\n# line break) untouched.\n\n$d = $word.ActiveDocument\n\n# Locate the start of the long text (the beginning of the offending run).\n$startRange = $d.Content.Duplicate\n$startFind = $startRange.Find\n$startFind.ClearFormatting()\n$startFind.Text = \"aqlFeatureAccess(org.eclipse.emf.ecore.EObject,java.lang.String) with arguments\"\n$startFound = $startFind.Execute()\n\n# Locate the end of the long text: the last line of the stack trace plus\n# the trailing newline character that precedes the run's closing tag.\n$endRange = $d.Content.Duplicate\n$endFind = $endRange.Find\n$endFind.ClearFormatting()\n$endFind.Text = \"Main.main(Main.java:1472)`n\"\n$endFound = $endFind.Execute()\n\nif ($startFound -and $endFound) {\n    # Build a single range spanning the whole stack trace and swap its\n    # text for the short replacement (formatting of the host run, i.e.\n    # bold + orange color, is preserved automatically).\n    $full = $d.Range($startRange.Start, $endRange.End)\n    $full.Text = \"Feature name not found in EClass EObject\"\n}\n"}
